$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2..17 (A..T) reflecting updated natmi LR-pair stats (3 samples per group)
$data = @(
    @("ECs", "Anxa1", "Fpr1", "ECs", 3, 1, 50.24752933333333, 150.742588, 0.1294604271951564, 0.1294604271951564, 2, 0.6666666666666666, 13.69437566666667, 41.083127, 0.6464706592402053, 0.6464706592402053, 688.1085430125195, 6192.976887112676, 0.08369236771437133, 0.08369236771437134),
    @("ECs", "Anxa1", "Fpr1", "FAPs", 3, 1, 50.24752933333333, 150.742588, 0.1294604271951564, 0.1294604271951564, 2, 0.6666666666666666, 2.346919, 7.040756999999999, 0.1107910510156661, 0.1107910510156661, 117.9268812954573, 1061.341931659116, 0.01434305679388849, 0.0143430567938885),
    @("ECs", "Anxa1", "Fpr1", "M2", 3, 1, 50.24752933333333, 150.742588, 0.1294604271951564, 0.1294604271951564, 3, 1, 3.811587333333333, 11.434762, 0.17993367759944, 0.17993367759944, 191.5228463382284, 1723.705617044056, 0.02329429076881903, 0.02329429076881904),
    @("ECs", "Anxa1", "Fpr1", "sCs", 3, 1, 50.24752933333333, 150.742588, 0.1294604271951564, 0.1294604271951564, 3, 1, 1.330408333333333, 3.991225, 0.06280461214468869, 0.06280461214468869, 66.84973175447776, 601.6475857902999, 0.008130711918077502, 0.008130711918077503),
    @("FAPs", "Anxa1", "Fpr1", "ECs", 3, 1, 151.42276, 454.26828, 0.3901337131747328, 0.3901337131747328, 2, 0.6666666666666666, 13.69437566666667, 41.083127, 0.6464706592402053, 0.6464706592402053, 2073.640159923507, 18662.76143931156, 0.2522099987478987, 0.2522099987478987),
    @("FAPs", "Anxa1", "Fpr1", "FAPs", 3, 1, 151.42276, 454.26828, 0.3901337131747328, 0.3901337131747328, 2, 0.6666666666666666, 2.346919, 7.040756999999999, 0.1107910510156661, 0.1107910510156661, 355.37695247644, 3198.39257228796, 0.04322332411927306, 0.04322332411927306),
    @("FAPs", "Anxa1", "Fpr1", "M2", 3, 1, 151.42276, 454.26828, 0.3901337131747328, 0.3901337131747328, 3, 1, 3.811587333333333, 11.434762, 0.17993367759944, 0.17993367759944, 577.1610739943734, 5194.44966594936, 0.07019819376705476, 0.07019819376705475),
    @("FAPs", "Anxa1", "Fpr1", "sCs", 3, 1, 151.42276, 454.26828, 0.3901337131747328, 0.3901337131747328, 3, 1, 1.330408333333333, 3.991225, 0.06280461214468869, 0.06280461214468869, 201.4541017603333, 1813.086915843, 0.02450219654050631, 0.02450219654050631),
    @("M2", "Anxa1", "Fpr1", "ECs", 3, 1, 85.147481, 255.442443, 0.2193785328573129, 0.2193785328573129, 2, 0.6666666666666666, 13.69437566666667, 41.083127, 0.6464706592402053, 0.6464706592402053, 1166.041591884362, 10494.37432695926, 0.1418217847594161, 0.1418217847594161),
    @("M2", "Anxa1", "Fpr1", "FAPs", 3, 1, 85.147481, 255.442443, 0.2193785328573129, 0.2193785328573129, 2, 0.6666666666666666, 2.346919, 7.040756999999999, 0.1107910510156661, 0.1107910510156661, 199.834240961039, 1798.508168649351, 0.02430517822553653, 0.02430517822553653),
    @("M2", "Anxa1", "Fpr1", "M2", 3, 1, 85.147481, 255.442443, 0.2193785328573129, 0.2193785328573129, 3, 1, 3.811587333333333, 11.434762, 0.17993367759944, 0.17993367759944, 324.5470600448406, 2920.923540403566, 0.03947358620338589, 0.03947358620338588),
    @("M2", "Anxa1", "Fpr1", "sCs", 3, 1, 85.147481, 255.442443, 0.2193785328573129, 0.2193785328573129, 3, 1, 1.330408333333333, 3.991225, 0.06280461214468869, 0.06280461214468869, 113.2809182847417, 1019.528264562675, 0.01377798366897438, 0.01377798366897438),
    @("sCs", "Anxa1", "Fpr1", "ECs", 3, 1, 101.3126446666667, 303.937934, 0.2610273267727979, 0.2610273267727979, 2, 0.6666666666666666, 13.69437566666667, 41.083127, 0.6464706592402053, 0.6464706592402053, 1387.413415848847, 12486.72074263962, 0.1687465080185192, 0.1687465080185192),
    @("sCs", "Anxa1", "Fpr1", "FAPs", 3, 1, 101.3126446666667, 303.937934, 0.2610273267727979, 0.2610273267727979, 2, 0.6666666666666666, 2.346919, 7.040756999999999, 0.1107910510156661, 0.1107910510156661, 237.7725707084487, 2139.953136376038, 0.028919491876968, 0.028919491876968),
    @("sCs", "Anxa1", "Fpr1", "M2", 3, 1, 101.3126446666667, 303.937934, 0.2610273267727979, 0.2610273267727979, 3, 1, 3.811587333333333, 11.434762, 0.17993367759944, 0.17993367759944, 386.1619931179676, 3475.457938061708, 0.0469676068601803, 0.04696760686018028),
    @("sCs", "Anxa1", "Fpr1", "sCs", 3, 1, 101.3126446666667, 303.937934, 0.2610273267727979, 0.2610273267727979, 3, 1, 1.330408333333333, 3.991225, 0.06280461214468869, 0.06280461214468869, 134.7871867365722, 1213.08468062915, 0.01639372001713049, 0.01639372001713049)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowVals[$c]
    }
}
